$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44160
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("P3").Value = 362

# Row 4
$ws.Range("D4").Value = 44167
$ws.Range("H4").Value = "Española"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("N4").Value = "`$/caja 30 unidades"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 450
$ws.Range("Q4").Value = 30

# Row 5
$ws.Range("D5").Value = 44405
$ws.Range("H5").Value = "Madrigal"
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 21000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 21500
$ws.Range("P5").Value = 538

# Row 6
$ws.Range("D6").Value = 44419
$ws.Range("J6").Value = 150
$ws.Range("N6").Value = "`$/caja 50 unidades"
$ws.Range("P6").Value = 430
$ws.Range("Q6").Value = 50

# Row 7
$ws.Range("D7").Value = 44370
$ws.Range("H7").Value = "Argentina(o)"
$ws.Range("J7").Value = 140
$ws.Range("M7").Value = 20429
$ws.Range("N7").Value = "`$/caja 50 unidades"
$ws.Range("P7").Value = 409
$ws.Range("Q7").Value = 50

# Row 8
$ws.Range("D8").Value = 44370
$ws.Range("H8").Value = "Madrigal"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 22000
$ws.Range("L8").Value = 23000
$ws.Range("M8").Value = 22500
$ws.Range("P8").Value = 562

# Row 9
$ws.Range("D9").Value = 44412
$ws.Range("H9").Value = "Symphony"
$ws.Range("J9").Value = 240
$ws.Range("K9").Value = 21000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 21500
$ws.Range("P9").Value = 538

# Row 10
$ws.Range("D10").Value = 44356
$ws.Range("H10").Value = "Argentina(o)"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("N10").Value = "`$/caja 50 unidades"
$ws.Range("P10").Value = 390
$ws.Range("Q10").Value = 50

# Row 11
$ws.Range("D11").Value = 44426
$ws.Range("H11").Value = "Madrigal"
$ws.Range("K11").Value = 19000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19500
$ws.Range("N11").Value = "`$/caja 40 unidades"
$ws.Range("P11").Value = 488
$ws.Range("Q11").Value = 40

# Row 12
$ws.Range("D12").Value = 44391
$ws.Range("J12").Value = 140

# Row 13
$ws.Range("D13").Value = 44435
$ws.Range("J13").Value = 160

# Row 14
$ws.Range("D14").Value = 44377
$ws.Range("H14").Value = "Madrigal"
$ws.Range("J14").Value = 150
$ws.Range("M14").Value = 20333
$ws.Range("N14").Value = "`$/caja 40 unidades"
$ws.Range("P14").Value = 508
$ws.Range("Q14").Value = 40

# Row 15
$ws.Range("D15").Value = 44377
$ws.Range("H15").Value = "Symphony"
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 21000
$ws.Range("L15").Value = 22000
$ws.Range("M15").Value = 21500
$ws.Range("P15").Value = 538

# Row 16
$ws.Range("D16").Value = 44433
$ws.Range("H16").Value = "Madrigal"
$ws.Range("K16").Value = 19000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 19500
$ws.Range("N16").Value = "`$/caja 40 unidades"
$ws.Range("O16").Value = "Región de Coquimbo"
$ws.Range("P16").Value = 488
$ws.Range("Q16").Value = 40

# Row 18
$ws.Range("D18").Value = 44384
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 21000
$ws.Range("L18").Value = 22000
$ws.Range("M18").Value = 21500
$ws.Range("P18").Value = 538

# Row 19
$ws.Range("D19").Value = 44384
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 30
$ws.Range("M19").Value = 19333
$ws.Range("N19").Value = "`$/caja 50 unidades"
$ws.Range("P19").Value = 387
$ws.Range("Q19").Value = 50

# Row 20
$ws.Range("D20").Value = 44384
$ws.Range("H20").Value = "Symphony"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 20000
$ws.Range("L20").Value = 21000
$ws.Range("M20").Value = 20400
$ws.Range("N20").Value = "`$/caja 40 unidades"
$ws.Range("P20").Value = 510
$ws.Range("Q20").Value = 40

# Row 21
$ws.Range("D21").Value = 44363
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19500
$ws.Range("P21").Value = 488
